$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")
Write-Host $ws.Range("A1").Value
